# CA 4.0 files test
# Applies the BDSBaPCF workbook update: refreshed plant-type list (incl. CCS /
# hydrogen / SMR variants), new "coal set to 0" note on the About sheet, and
# header/label rewording.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# About sheet: append the three new note lines (rows 24-26)
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A24").Value2 = "For the United States, we have set coal to 0 as of version 3.4. This reflects"
$wsAbout.Range("A25").Value2 = "the fact that certain air quality / environmental restrictions, as well as current"
$wsAbout.Range("A26").Value2 = "supply chain logistics, limit the amount the coal dispatches annually. "

# ---------------------------------------------------------------------------
# BDSBaPCF sheet: rebuild the data table with the new plant list / order
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("BDSBaPCF")

# Clear the old table (values + formulas) before laying out the new one.
$wsData.Range("A1:B17").ClearContents() | Out-Null

# Header row
$wsData.Range("A1").Value2 = "Unit: dimensionless (Boolean)"
$wsData.Range("B1").Value2 = "Do Suppliers Bid at Peak Capacity Factors"
$wsData.Rows.Item(1).RowHeight = 30

# Plant types bidding at peak capacity factor (directly entered)
$wsData.Range("A2").Value2 = "hard coal"
$wsData.Range("B2").Value2 = 0

$wsData.Range("A3").Value2 = "natural gas steam turbine"
$wsData.Range("B3").Value2 = 1

$wsData.Range("A4").Value2 = "natural gas combined cycle"
$wsData.Range("B4").Value2 = 1

$wsData.Range("A5").Value2 = "nuclear"
$wsData.Range("B5").Value2 = 0

$wsData.Range("A6").Value2 = "hydro"
$wsData.Range("B6").Value2 = 0

$wsData.Range("A7").Value2 = "onshore wind"
$wsData.Range("B7").Value2 = 0

$wsData.Range("A8").Value2 = "solar PV"
$wsData.Range("B8").Value2 = 0

$wsData.Range("A9").Value2 = "solar thermal"
$wsData.Range("B9").Value2 = 0

$wsData.Range("A10").Value2 = "biomass"
$wsData.Range("B10").Value2 = 1

$wsData.Range("A11").Value2 = "geothermal"
$wsData.Range("B11").Value2 = 0

$wsData.Range("A12").Value2 = "petroleum"
$wsData.Range("B12").Value2 = 0

$wsData.Range("A13").Value2 = "natural gas peaker"
$wsData.Range("B13").Value2 = 0

# Plant types that mirror another plant type's bid (formulas)
$wsData.Range("A14").Value2 = "lignite"
$wsData.Range("B14").Formula = "=B2"

$wsData.Range("A15").Value2 = "offshore wind"
$wsData.Range("B15").Formula = "=B7"

$wsData.Range("A16").Value2 = "crude oil"
$wsData.Range("B16").Formula = "=B12"

$wsData.Range("A17").Value2 = "heavy or residual fuel oil"
$wsData.Range("B17").Formula = "=B12"

$wsData.Range("A18").Value2 = "municipal solid waste"
$wsData.Range("B18").Formula = "=B10"

# New CCS / advanced-nuclear / hydrogen plant types (highlighted, integer format)
$wsData.Range("A19").Value2 = "hard coal w CCS"
$wsData.Range("B19").Value2 = 0

$wsData.Range("A20").Value2 = "natural gas combined cycle w CCS"
$wsData.Range("B20").Value2 = 0

$wsData.Range("A21").Value2 = "biomass w CCS"
$wsData.Range("B21").Value2 = 0

$wsData.Range("A22").Value2 = "lignite w CCS"
$wsData.Range("B22").Value2 = 0

$wsData.Range("A23").Value2 = "small modular reactor"
$wsData.Range("B23").Value2 = 0

$wsData.Range("A24").Value2 = "hydrogen combustion turbine"
$wsData.Range("B24").Value2 = 0

$wsData.Range("A25").Value2 = "hydrogen combined cycle"
$wsData.Range("B25").Value2 = 0

# Highlight fill + integer number format on the new boolean column cells.
$rCcs = $wsData.Range("B19:B25")
$rCcs.NumberFormat = "0"
$rCcs.Interior.ThemeColor = 8

# The two hydrogen rows get a slightly different (black, vertically centered) label font.
$rHyd = $wsData.Range("A24:A25")
$rHyd.Font.Color = 0
$rHyd.VerticalAlignment = -4108

# Unit header label is italicized.
$wsData.Range("A1").Font.Italic = $true

$wsData.Columns.Item(1).ColumnWidth = 29.29

# View state: BDSBaPCF becomes the active/selected sheet, with D32 selected.
$wsData.Activate() | Out-Null
$wsData.Range("D32").Select() | Out-Null
